# Remove the RG (ID card) reference from the certificate text, leaving
# only the CPF reference, and fold a couple of adjacent same-formatted
# runs back together the way PowerPoint does when text is retyped.
#
# Target shape: Slide 1, "Rectangle 5" - the body paragraph that reads
#   "Certificamos que {{NOME}}, portador do RG nº {{RG}} e CPF nº {{CPF}}, ..."
# becomes
#   "Certificamos que {{NOME}}, portador do CPF nº {{CPF}}, ..."
# and, lower in the same shape, the signature-date line
#   "São Carlos" + ", " + "{{DATA}}"  ->  "São Carlos, " + "{{DATA}}"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item("Rectangle 5")
$tf = $shp.TextFrame
$tr = $tf.TextRange

# Locate the (0-based) offsets of the relevant substrings in the
# shape's current text so the edit doesn't depend on brittle, hand
# counted character positions.
$full = $tr.Text
$iPortador = $full.IndexOf("portador do RG")
$iRG       = $full.IndexOf("{{RG}} ")
$iECPF     = $full.IndexOf("e CPF n")
$iConcl    = $full.IndexOf(", concluiu com aproveitamento o curso da ")
$iCarlos   = $full.IndexOf("Carlos") - 4

# Apply the edits from the rightmost position back to the leftmost so
# that earlier, still-to-be-used offsets are never invalidated by a
# preceding (length changing) replacement.

# "São Carlos" + ", " -> a single run "São Carlos, "
$rCarlos = $tr.Characters($iCarlos + 1, 12)
$rCarlos.Text = "São Carlos, "

# ", " + "concluiu com aproveitamento o curso da " -> a single run
$rConcl = $tr.Characters($iConcl + 1, 41)
$rConcl.Text = ", concluiu com aproveitamento o curso da "

# "e CPF nº " -> "nº "
$rECPF = $tr.Characters($iECPF + 1, 9)
$rECPF.Text = "nº "

# second half of the old "{{RG}} " run -> "CPF "
$rCPF = $tr.Characters($iRG + 1 + 3, 4)
$rCPF.Text = "CPF "

# first half of the old "{{RG}} " run -> "do "
$rDo = $tr.Characters($iRG + 1, 3)
$rDo.Text = "do "

# "portador do RG nº " -> "portador "
$rPortador = $tr.Characters($iPortador + 1, 18)
$rPortador.Text = "portador "
